# Apply the "Depencency inversion example added" edit to the SOLID
# principles document:
#   1. Collapse the "Single Responsibility..." paragraph's three runs
#      (split around a gramStart/gramEnd proofErr pair caused by the
#      double space in "have  a") into a single run with the same text.
#   2. Append three new paragraphs after it: the Open/Closed principle,
#      the Liskov Substitution principle, and the start of the Interface
#      segregation paragraph (which carries a lastRenderedPageBreak
#      marker, as Word stamps on text that begins a new page).

$d = $word.ActiveDocument

$enDash = [char]0x2013

# --- 1. Merge the Single Responsibility Principle paragraph into one run ---
# Replacing the whole phrase (which spans all three existing runs) with
# itself via Find/Replace collapses it back down to a single run and drops
# the proofErr gramStart/gramEnd markers that bracketed "have  a".
$srpText = "Single Responsibility Principle " + $enDash + " a class should have  a single reason to change."
$found = $d.Content.Find.Execute($srpText, $true, $false, $false, $false, $false, `
                                  $true, 1, $false, $srpText, 2)
if (-not $found) {
    throw "Could not locate the Single Responsibility Principle paragraph text"
}

# --- 2. Append the three new paragraphs after it ---

function Insert-ParagraphXml([string]$innerParagraphXml) {
    $lastPara = $d.Paragraphs.Last
    $insertionPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
           '<w:body>' + $innerParagraphXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $null = $insertionPoint.InsertXML($pkg)
}

Insert-ParagraphXml('<w:p><w:r><w:t>Open close- your app should be open for inheritance, but closed for modification</w:t></w:r></w:p>')

Insert-ParagraphXml('<w:p><w:r><w:t>Liskov Substitution ' + $enDash + ' any function that takes as argument a Base class object, should work perfectly fine with a derived class object.</w:t></w:r></w:p>')

Insert-ParagraphXml('<w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Interface segregation : </w:t></w:r></w:p>')
